$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"
# Card number is a 16-digit string that must stay text (not be coerced to a
# number) -- a leading apostrophe forces Excel to treat it as text while
# keeping the cell's existing "General" number format (matches how this
# value is stored as a literal string in the source file).
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 21.07.2025"

# Transaction rows
$ws.Range("B6").Value = "22.07."
$ws.Range("C6").Value = "23.07."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-66466022"
$ws.Range("E6").Value = "55,49-"

$ws.Range("B7").Value = "24.07."
$ws.Range("C7").Value = "25.07."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,19-"

$ws.Range("B8").Value = "27.07."
$ws.Range("C8").Value = "28.07."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 63775419"
$ws.Range("E8").Value = "40,44-"

$ws.Range("B9").Value = "30.07."
$ws.Range("C9").Value = "31.07."
$ws.Range("D9").Value = "PAYPAL RDCOGF"
$ws.Range("E9").Value = "88,74-"

$ws.Range("B10").Value = "03.08."
$ws.Range("C10").Value = "04.08."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 92254201"
$ws.Range("E10").Value = "84,53-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 06.08.2025"
$ws.Range("E12").Value = "294,39-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.08.2025"
